$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Formula = "=3.75-0.0701941"
$ws.Range("D8").Formula = "=3.5 - 0.0548058971"
$ws.Range("B9").Value = 0.75

$ws.Range("D9").Select()

$wb.Application.CalculateFull()
